$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Reposition the two floating drawings (picture + caption text box)
#    Shapes collection is in document order:
#      Item(1) -> "Picture 1"  (wp14:anchorId="7B095336")
#      Item(2) -> "Text Box 2" (wp14:anchorId="2B24792C")
#    Shape.Left/Top are expressed in points; 1 pt = 12700 EMU.
# ------------------------------------------------------------------

# Picture 1: posOffset 0,572135 (EMU) -> 4457700,343535 (EMU)
$picture = $d.Shapes.Item(1)
$picture.Left = 4457700 / 12700   # 351 pt
$picture.Top  = 343535 / 12700    # 27.05 pt

# Text Box 2: posOffset 1714500,377825 (EMU) -> -342900,263525 (EMU)
$textbox = $d.Shapes.Item(2)
$textbox.Left = -342900 / 12700   # -27 pt
$textbox.Top  = 263525 / 12700    # 20.75 pt

# ------------------------------------------------------------------
# 2) Move the (empty) "_GoBack" bookmark from the paragraph that hosts
#    the caption text box to the end of the preceding (blank) paragraph.
#    Re-adding a bookmark with the same name relocates it (bookmark
#    names are unique within a document).
# ------------------------------------------------------------------

$paragraphs = $d.Paragraphs
$target = $null
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $para = $paragraphs.Item($i)
    if ($para.Range.Text -eq "`r") {
        $target = $para
    }
}
if ($target -ne $null) {
    $d.Bookmarks.Add("_GoBack", $target.Range)
}
